# Resource Flow 29 November
# Update the logged UserName / WorkGroup resource-flow entry in row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "UserName1542882306701"
$ws.Range("C2").Value = "WorkGroup1542882450361"
